$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 ("SOURCES OF FINANCE") table: switch the applied table style.
#    {8944E028-05A9-4674-B7E6-53E175123D58} -> {4C469218-8B9D-4220-80FB-DEEC781A6B4D}
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{4C469218-8B9D-4220-80FB-DEEC781A6B4D}")
    }
}

# ---------------------------------------------------------------------------
# 2) Deck theme colour scheme: move from the "Integral" palette to the
#    stock "Office Theme" palette (dk1/lt1 are already shared by both).
# ---------------------------------------------------------------------------
$theme = $p.Designs.Item(1).SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
